$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.782436333333333
$ws.Range("H2").Value = 5.347308999999999
$ws.Range("I2").Value = 0.00914036392049929
$ws.Range("J2").Value = 0.009140363920499292
$ws.Range("M2").Value = 7.579746333333333
$ws.Range("N2").Value = 22.739239
$ws.Range("O2").Value = 0.0686314777863378
$ws.Range("P2").Value = 0.0686314777863378
$ws.Range("Q2").Value = 13.51041526198344
$ws.Range("R2").Value = 121.593737357851
$ws.Range("S2").Value = 0.0006273166833687904
$ws.Range("T2").Value = 0.0006273166833687906
$ws.Range("G3").Value = 1.782436333333333
$ws.Range("H3").Value = 5.347308999999999
$ws.Range("I3").Value = 0.00914036392049929
$ws.Range("J3").Value = 0.009140363920499292
$ws.Range("O3").Value = 0.0596740760116217
$ws.Range("P3").Value = 0.05967407601162171
$ws.Range("Q3").Value = 11.74711041196122
$ws.Range("R3").Value = 105.723993707651
$ws.Range("S3").Value = 0.0005454427713657592
$ws.Range("T3").Value = 0.0005454427713657594
$ws.Range("G4").Value = 1.782436333333333
$ws.Range("H4").Value = 5.347308999999999
$ws.Range("I4").Value = 0.00914036392049929
$ws.Range("J4").Value = 0.009140363920499292
$ws.Range("M4").Value = 96.08192699999999
$ws.Range("N4").Value = 288.245781
$ws.Range("O4").Value = 0.8699822327258658
$ws.Range("P4").Value = 0.8699822327258659
$ws.Range("Q4").Value = 171.259917661481
$ws.Range("R4").Value = 1541.339258953329
$ws.Range("S4").Value = 0.00795195421148292
$ws.Range("T4").Value = 0.007951954211482924
$ws.Range("G5").Value = 1.782436333333333
$ws.Range("H5").Value = 5.347308999999999
$ws.Range("I5").Value = 0.00914036392049929
$ws.Range("J5").Value = 0.009140363920499292
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.189099
$ws.Range("N5").Value = 0.5672970000000001
$ws.Range("O5").Value = 0.001712213476174646
$ws.Range("P5").Value = 0.001712213476174646
$ws.Range("Q5").Value = 0.337056928197
$ws.Range("R5").Value = 3.033512353773
$ws.Range("S5").Value = [double]"1.56502542818194E-05"
$ws.Range("T5").Value = [double]"1.56502542818194E-05"
$ws.Range("H6").Value = 564.692825
$ws.Range("I6").Value = 0.965251479537618
$ws.Range("J6").Value = 0.965251479537618
$ws.Range("M6").Value = 7.579746333333333
$ws.Range("N6").Value = 22.739239
$ws.Range("O6").Value = 0.0686314777863378
$ws.Range("P6").Value = 0.0686314777863378
$ws.Range("Q6").Value = 1426.742789917797
$ws.Range("R6").Value = 12840.68510926018
$ws.Range("S6").Value = 0.06624663547611573
$ws.Range("T6").Value = 0.06624663547611573
$ws.Range("H7").Value = 564.692825
$ws.Range("I7").Value = 0.965251479537618
$ws.Range("J7").Value = 0.965251479537618
$ws.Range("O7").Value = 0.0596740760116217
$ws.Range("P7").Value = 0.05967407601162171
$ws.Range("S7").Value = 0.05760049016025813
$ws.Range("T7").Value = 0.05760049016025814
$ws.Range("H8").Value = 564.692825
$ws.Range("I8").Value = 0.965251479537618
$ws.Range("J8").Value = 0.965251479537618
$ws.Range("M8").Value = 96.08192699999999
$ws.Range("N8").Value = 288.245781
$ws.Range("O8").Value = 0.8699822327258658
$ws.Range("P8").Value = 0.8699822327258659
$ws.Range("Q8").Value = 18085.59159635792
$ws.Range("R8").Value = 162770.3243672213
$ws.Range("S8").Value = 0.8397516373100823
$ws.Range("T8").Value = 0.8397516373100824
$ws.Range("H9").Value = 564.692825
$ws.Range("I9").Value = 0.965251479537618
$ws.Range("J9").Value = 0.965251479537618
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.189099
$ws.Range("N9").Value = 0.5672970000000001
$ws.Range("O9").Value = 0.001712213476174646
$ws.Range("P9").Value = 0.001712213476174646
$ws.Range("Q9").Value = 35.59428283822501
$ws.Range("R9").Value = 320.348545544025
$ws.Range("S9").Value = 0.001652716591161825
$ws.Range("T9").Value = 0.001652716591161825
$ws.Range("G10").Value = 4.870778333333333
$ws.Range("H10").Value = 14.612335
$ws.Range("I10").Value = 0.02497743437460768
$ws.Range("J10").Value = 0.02497743437460768
$ws.Range("M10").Value = 7.579746333333333
$ws.Range("N10").Value = 22.739239
$ws.Range("O10").Value = 0.0686314777863378
$ws.Range("P10").Value = 0.0686314777863378
$ws.Range("Q10").Value = 36.91926421256277
$ws.Range("R10").Value = 332.273377913065
$ws.Range("S10").Value = 0.001714238232440597
$ws.Range("T10").Value = 0.001714238232440597
$ws.Range("G11").Value = 4.870778333333333
$ws.Range("H11").Value = 14.612335
$ws.Range("I11").Value = 0.02497743437460768
$ws.Range("J11").Value = 0.02497743437460768
$ws.Range("O11").Value = 0.0596740760116217
$ws.Range("P11").Value = 0.05967407601162171
$ws.Range("Q11").Value = 32.10076556667389
$ws.Range("R11").Value = 288.906890100065
$ws.Range("S11").Value = 0.001490505317445631
$ws.Range("T11").Value = 0.001490505317445631
$ws.Range("G12").Value = 4.870778333333333
$ws.Range("H12").Value = 14.612335
$ws.Range("I12").Value = 0.02497743437460768
$ws.Range("J12").Value = 0.02497743437460768
$ws.Range("M12").Value = 96.08192699999999
$ws.Range("N12").Value = 288.245781
$ws.Range("O12").Value = 0.8699822327258658
$ws.Range("P12").Value = 0.8699822327258659
$ws.Range("Q12").Value = 467.993768256515
$ws.Range("R12").Value = 4211.943914308634
$ws.Range("S12").Value = 0.02172992412498498
$ws.Range("T12").Value = 0.02172992412498498
$ws.Range("G13").Value = 4.870778333333333
$ws.Range("H13").Value = 14.612335
$ws.Range("I13").Value = 0.02497743437460768
$ws.Range("J13").Value = 0.02497743437460768
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.189099
$ws.Range("N13").Value = 0.5672970000000001
$ws.Range("O13").Value = 0.001712213476174646
$ws.Range("P13").Value = 0.001712213476174646
$ws.Range("Q13").Value = 0.921059312055
$ws.Range("R13").Value = 8.289533808494999
$ws.Range("S13").Value = [double]"4.27666997364711E-05"
$ws.Range("T13").Value = [double]"4.27666997364711E-05"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1229953333333333
$ws.Range("H14").Value = 0.368986
$ws.Range("I14").Value = 0.0006307221672750447
$ws.Range("J14").Value = 0.0006307221672750447
$ws.Range("M14").Value = 7.579746333333333
$ws.Range("N14").Value = 22.739239
$ws.Range("O14").Value = 0.0686314777863378
$ws.Range("P14").Value = 0.0686314777863378
$ws.Range("Q14").Value = 0.9322734268504446
$ws.Range("R14").Value = 8.390460841654001
$ws.Range("S14").Value = [double]"4.328739441268806E-05"
$ws.Range("T14").Value = [double]"4.328739441268806E-05"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1229953333333333
$ws.Range("H15").Value = 0.368986
$ws.Range("I15").Value = 0.0006307221672750447
$ws.Range("J15").Value = 0.0006307221672750447
$ws.Range("O15").Value = 0.0596740760116217
$ws.Range("P15").Value = 0.05967407601162171
$ws.Range("Q15").Value = 0.8105982434282223
$ws.Range("R15").Value = 7.295384190854001
$ws.Range("S15").Value = [double]"3.76377625521858E-05"
$ws.Range("T15").Value = [double]"3.76377625521858E-05"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1229953333333333
$ws.Range("H16").Value = 0.368986
$ws.Range("I16").Value = 0.0006307221672750447
$ws.Range("J16").Value = 0.0006307221672750447
$ws.Range("M16").Value = 96.08192699999999
$ws.Range("N16").Value = 288.245781
$ws.Range("O16").Value = 0.8699822327258658
$ws.Range("P16").Value = 0.8699822327258659
$ws.Range("Q16").Value = 11.817628638674
$ws.Range("R16").Value = 106.358657748066
$ws.Range("S16").Value = 0.0005487170793156404
$ws.Range("T16").Value = 0.0005487170793156405
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.1229953333333333
$ws.Range("H17").Value = 0.368986
$ws.Range("I17").Value = 0.0006307221672750447
$ws.Range("J17").Value = 0.0006307221672750447
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.189099
$ws.Range("N17").Value = 0.5672970000000001
$ws.Range("O17").Value = 0.001712213476174646
$ws.Range("P17").Value = 0.001712213476174646
$ws.Range("Q17").Value = 0.023258294538
$ws.Range("R17").Value = 0.209324650842
$ws.Range("S17").Value = [double]"1.079930994530411E-06"
$ws.Range("T17").Value = [double]"1.079930994530411E-06"
